$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Tercer atributo"
$ws.Range("B6").Value = "ojnfwjf"
$ws.Range("C6").Value = "owirjfw"
$ws.Range("D6").Value = "iwejrgpwergf"
